$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily records to append (aggiornamento fino a 20/09/2021)
# columns: date-serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$data = @(
    @(44449, 1, 2, 50.8646998982706),
    @(44450, 0, 2, 50.8646998982706),
    @(44451, 2, 4, 101.7293997965412),
    @(44452, 1, 5, 127.1617497456765),
    @(44453, 5, 10, 254.323499491353),
    @(44454, 0, 10, 254.323499491353),
    @(44455, 0, 9, 228.8911495422177),
    @(44456, 2, 10, 254.323499491353),
    @(44457, 0, 10, 254.323499491353),
    @(44458, 0, 8, 203.4587995930824),
    @(44459, 1, 8, 203.4587995930824)
)

$lastRow = 374
$startRow = $lastRow + 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    # Carry the formatting (date style etc.) from the previous row down to the new one.
    $src = $ws.Range("A" + $lastRow + ":D" + $lastRow)
    $dst = $ws.Range("A" + $r + ":D" + $r)
    $src.Copy($dst)

    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
